$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new data row (row 8) carried over from the source feed.
$ws.Range("A8").Value = 42612.890752314815
$ws.Range("B8").Value = -24
$ws.Range("C8").Value = 47
$ws.Range("D8").Value = 51
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 99
$ws.Range("G8").Value = 13807
$ws.Range("H8").Value = 10794
$ws.Range("I8").Value = 580
$ws.Range("J8").Value = 86
$ws.Range("K8").Value = 93
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 7
$ws.Range("N8").Value = "Named"
